$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "I" column (header "sqrt" + per-row SQRT(ABS(D..)) formulas) moves one
# column to the left, into "H" (which was empty before).
$ws.Range("H1").Value = $ws.Range("I1").Value2
$ws.Range("H1").NumberFormat = $ws.Range("I1").NumberFormat
$ws.Range("I1").Clear()

for ($r = 2; $r -le 14; $r++) {
    $ws.Range("H$r").Formula = $ws.Range("I$r").Formula
    $ws.Range("I$r").Clear()
}

# The "K" column (header "liest" + constant values) is removed entirely.
$ws.Range("K1:K14").Clear()

# Restore the cursor/selection to where the user last clicked.
[void]$ws.Range("I18").Select()
